$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate rows 18-20 (F:V): row18<-row20, row19<-row18, row20<-row19 ---
# Row 18 gets old Row 20 data (F:V)
$ws.Cells.Item(18, 6).Value = 'Buducnost'
$ws.Cells.Item(18, 7).Value = 1
$ws.Cells.Item(18, 8).Value = 'Jezero'
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(18, 10).Value = 1.48
$ws.Cells.Item(18, 11).Value = '12/08/2023 09:12'
$ws.Cells.Item(18, 12).Value = 1.47
$ws.Cells.Item(18, 13).Value = '13/08/2023 19:59'
$ws.Cells.Item(18, 14).Value = 3.83
$ws.Cells.Item(18, 15).Value = '12/08/2023 09:12'
$ws.Cells.Item(18, 16).Value = 3.91
$ws.Cells.Item(18, 17).Value = '13/08/2023 19:59'
$ws.Cells.Item(18, 18).Value = 5.58
$ws.Cells.Item(18, 19).Value = '12/08/2023 09:12'
$ws.Cells.Item(18, 20).Value = 7.31
$ws.Cells.Item(18, 21).Value = '13/08/2023 19:59'
$ws.Cells.Item(18, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/buducnost-jezero/nXQrD1KH/'

# Row 19 gets old Row 18 data (F:V)
$ws.Cells.Item(19, 6).Value = 'Mladost DG'
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 'Decic'
$ws.Cells.Item(19, 9).Value = 2
$ws.Cells.Item(19, 10).Value = 2.64
$ws.Cells.Item(19, 11).Value = '13/08/2023 10:35'
$ws.Cells.Item(19, 12).Value = 3.16
$ws.Cells.Item(19, 13).Value = '13/08/2023 19:49'
$ws.Cells.Item(19, 14).Value = 2.99
$ws.Cells.Item(19, 15).Value = '13/08/2023 10:35'
$ws.Cells.Item(19, 16).Value = 2.89
$ws.Cells.Item(19, 17).Value = '13/08/2023 19:39'
$ws.Cells.Item(19, 18).Value = 2.75
$ws.Cells.Item(19, 19).Value = '13/08/2023 10:35'
$ws.Cells.Item(19, 20).Value = 2.36
$ws.Cells.Item(19, 21).Value = '13/08/2023 19:49'
$ws.Cells.Item(19, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mladost-dg-decic/8nFiBukU/'

# Row 20 gets old Row 19 data (F:V)
$ws.Cells.Item(20, 6).Value = 'Jedinstvo'
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 'Sutjeska'
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 5.84
$ws.Cells.Item(20, 11).Value = '12/08/2023 17:12'
$ws.Cells.Item(20, 12).Value = 4.16
$ws.Cells.Item(20, 13).Value = '13/08/2023 19:20'
$ws.Cells.Item(20, 14).Value = 3.85
$ws.Cells.Item(20, 15).Value = '12/08/2023 17:12'
$ws.Cells.Item(20, 16).Value = 3.62
$ws.Cells.Item(20, 17).Value = '13/08/2023 19:20'
$ws.Cells.Item(20, 18).Value = 1.48
$ws.Cells.Item(20, 19).Value = '12/08/2023 17:12'
$ws.Cells.Item(20, 20).Value = 1.8
$ws.Cells.Item(20, 21).Value = '13/08/2023 19:20'
$ws.Cells.Item(20, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jedinstvo-sutjeska/YTMvEs5B/'

# --- Swap rows 24-25 (F:V) ---
# Row 24 gets old Row 25 data (F:V)
$ws.Cells.Item(24, 6).Value = 'Jezero'
$ws.Cells.Item(24, 7).Value = 2
$ws.Cells.Item(24, 8).Value = 'Petrovac'
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 2.46
$ws.Cells.Item(24, 11).Value = '19/08/2023 19:42'
$ws.Cells.Item(24, 12).Value = 2.65
$ws.Cells.Item(24, 13).Value = '20/08/2023 19:46'
$ws.Cells.Item(24, 14).Value = 2.86
$ws.Cells.Item(24, 15).Value = '19/08/2023 19:42'
$ws.Cells.Item(24, 16).Value = 2.99
$ws.Cells.Item(24, 17).Value = '20/08/2023 19:46'
$ws.Cells.Item(24, 18).Value = 2.8
$ws.Cells.Item(24, 19).Value = '19/08/2023 19:42'
$ws.Cells.Item(24, 20).Value = 2.78
$ws.Cells.Item(24, 21).Value = '20/08/2023 19:46'
$ws.Cells.Item(24, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jezero-petrovac/KQv6U2kh/'

# Row 25 gets old Row 24 data (F:V)
$ws.Cells.Item(25, 6).Value = 'Decic'
$ws.Cells.Item(25, 7).Value = 3
$ws.Cells.Item(25, 8).Value = 'Mornar Bar'
$ws.Cells.Item(25, 9).Value = 1
$ws.Cells.Item(25, 10).Value = 2.07
$ws.Cells.Item(25, 11).Value = '19/08/2023 19:42'
$ws.Cells.Item(25, 12).Value = 1.69
$ws.Cells.Item(25, 13).Value = '20/08/2023 18:23'
$ws.Cells.Item(25, 14).Value = 2.85
$ws.Cells.Item(25, 15).Value = '19/08/2023 19:42'
$ws.Cells.Item(25, 16).Value = 3.19
$ws.Cells.Item(25, 17).Value = '20/08/2023 18:23'
$ws.Cells.Item(25, 18).Value = 3.57
$ws.Cells.Item(25, 19).Value = '19/08/2023 19:42'
$ws.Cells.Item(25, 20).Value = 6.08
$ws.Cells.Item(25, 21).Value = '20/08/2023 18:23'
$ws.Cells.Item(25, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/decic-mornar-bar/bBjcWOKu/'

# --- Append new rows 41-45 ---
# Row 41: copy style from row 40 for columns A and E (bold index col / date col)
$ws.Range("A40").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("E40").Copy()
$ws.Range("E41").PasteSpecial(-4122)
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = 'montenegro'
$ws.Cells.Item(41, 3).Value = 'prva-crnogorska-liga'
$ws.Cells.Item(41, 4).Value = '2023-2024'
$ws.Cells.Item(41, 5).Value = 45192.66666666666
$ws.Cells.Item(41, 6).Value = 'Arsenal Tivat'
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 'Decic'
$ws.Cells.Item(41, 9).Value = 1
$ws.Cells.Item(41, 10).Value = 3.05
$ws.Cells.Item(41, 11).Value = '22/09/2023 03:13'
$ws.Cells.Item(41, 12).Value = 3.73
$ws.Cells.Item(41, 13).Value = '23/09/2023 15:51'
$ws.Cells.Item(41, 14).Value = 2.82
$ws.Cells.Item(41, 15).Value = '22/09/2023 03:13'
$ws.Cells.Item(41, 16).Value = 2.92
$ws.Cells.Item(41, 17).Value = '23/09/2023 15:51'
$ws.Cells.Item(41, 18).Value = 2.32
$ws.Cells.Item(41, 19).Value = '22/09/2023 03:13'
$ws.Cells.Item(41, 20).Value = 2.17
$ws.Cells.Item(41, 21).Value = '23/09/2023 15:51'
$ws.Cells.Item(41, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/arsenal-tivat-decic/2szPtXvF/'

# Row 42: copy style from row 40 for columns A and E (bold index col / date col)
$ws.Range("A40").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("E40").Copy()
$ws.Range("E42").PasteSpecial(-4122)
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = 'montenegro'
$ws.Cells.Item(42, 3).Value = 'prva-crnogorska-liga'
$ws.Cells.Item(42, 4).Value = '2023-2024'
$ws.Cells.Item(42, 5).Value = 45192.66666666666
$ws.Cells.Item(42, 6).Value = 'Sutjeska'
$ws.Cells.Item(42, 7).Value = 1
$ws.Cells.Item(42, 8).Value = 'Rudar'
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 1.34
$ws.Cells.Item(42, 11).Value = '22/09/2023 03:13'
$ws.Cells.Item(42, 12).Value = 1.44
$ws.Cells.Item(42, 13).Value = '23/09/2023 15:58'
$ws.Cells.Item(42, 14).Value = 4.52
$ws.Cells.Item(42, 15).Value = '22/09/2023 03:13'
$ws.Cells.Item(42, 16).Value = 4.67
$ws.Cells.Item(42, 17).Value = '23/09/2023 15:59'
$ws.Cells.Item(42, 18).Value = 7.27
$ws.Cells.Item(42, 19).Value = '22/09/2023 03:13'
$ws.Cells.Item(42, 20).Value = 6.14
$ws.Cells.Item(42, 21).Value = '23/09/2023 15:59'
$ws.Cells.Item(42, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/sutjeska-rudar/KIpKsDP8/'

# Row 43: copy style from row 40 for columns A and E (bold index col / date col)
$ws.Range("A40").Copy()
$ws.Range("A43").PasteSpecial(-4122)
$ws.Range("E40").Copy()
$ws.Range("E43").PasteSpecial(-4122)
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = 'montenegro'
$ws.Cells.Item(43, 3).Value = 'prva-crnogorska-liga'
$ws.Cells.Item(43, 4).Value = '2023-2024'
$ws.Cells.Item(43, 5).Value = 45192.70833333334
$ws.Cells.Item(43, 6).Value = 'Jezero'
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 'Mornar Bar'
$ws.Cells.Item(43, 9).Value = 1
$ws.Cells.Item(43, 10).Value = 2.21
$ws.Cells.Item(43, 11).Value = '22/09/2023 04:12'
$ws.Cells.Item(43, 12).Value = 2.45
$ws.Cells.Item(43, 13).Value = '23/09/2023 16:46'
$ws.Cells.Item(43, 14).Value = 2.78
$ws.Cells.Item(43, 15).Value = '22/09/2023 04:12'
$ws.Cells.Item(43, 16).Value = 2.75
$ws.Cells.Item(43, 17).Value = '23/09/2023 16:46'
$ws.Cells.Item(43, 18).Value = 3.32
$ws.Cells.Item(43, 19).Value = '22/09/2023 04:12'
$ws.Cells.Item(43, 20).Value = 3.33
$ws.Cells.Item(43, 21).Value = '23/09/2023 16:46'
$ws.Cells.Item(43, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jezero-mornar-bar/6FtGrgA2/'

# Row 44: copy style from row 40 for columns A and E (bold index col / date col)
$ws.Range("A40").Copy()
$ws.Range("A44").PasteSpecial(-4122)
$ws.Range("E40").Copy()
$ws.Range("E44").PasteSpecial(-4122)
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = 'montenegro'
$ws.Cells.Item(44, 3).Value = 'prva-crnogorska-liga'
$ws.Cells.Item(44, 4).Value = '2023-2024'
$ws.Cells.Item(44, 5).Value = 45192.75
$ws.Cells.Item(44, 6).Value = 'Jedinstvo'
$ws.Cells.Item(44, 7).Value = 2
$ws.Cells.Item(44, 8).Value = 'Mladost DG'
$ws.Cells.Item(44, 9).Value = 1
$ws.Cells.Item(44, 10).Value = 2.25
$ws.Cells.Item(44, 11).Value = '22/09/2023 05:13'
$ws.Cells.Item(44, 12).Value = 2.34
$ws.Cells.Item(44, 13).Value = '22/09/2023 11:41'
$ws.Cells.Item(44, 14).Value = 2.94
$ws.Cells.Item(44, 15).Value = '22/09/2023 05:13'
$ws.Cells.Item(44, 16).Value = 3.01
$ws.Cells.Item(44, 17).Value = '23/09/2023 16:04'
$ws.Cells.Item(44, 18).Value = 3.02
$ws.Cells.Item(44, 19).Value = '22/09/2023 05:13'
$ws.Cells.Item(44, 20).Value = 3.17
$ws.Cells.Item(44, 21).Value = '23/09/2023 10:26'
$ws.Cells.Item(44, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jedinstvo-mladost-dg/tYYSuifL/'

# Row 45: copy style from row 40 for columns A and E (bold index col / date col)
$ws.Range("A40").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("E40").Copy()
$ws.Range("E45").PasteSpecial(-4122)
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = 'montenegro'
$ws.Cells.Item(45, 3).Value = 'prva-crnogorska-liga'
$ws.Cells.Item(45, 4).Value = '2023-2024'
$ws.Cells.Item(45, 5).Value = 45192.79166666666
$ws.Cells.Item(45, 6).Value = 'Buducnost'
$ws.Cells.Item(45, 7).Value = 2
$ws.Cells.Item(45, 8).Value = 'Petrovac'
$ws.Cells.Item(45, 9).Value = 1
$ws.Cells.Item(45, 10).Value = 1.46
$ws.Cells.Item(45, 11).Value = '22/09/2023 15:12'
$ws.Cells.Item(45, 12).Value = 1.53
$ws.Cells.Item(45, 13).Value = '23/09/2023 18:00'
$ws.Cells.Item(45, 14).Value = 3.92
$ws.Cells.Item(45, 15).Value = '22/09/2023 15:12'
$ws.Cells.Item(45, 16).Value = 3.87
$ws.Cells.Item(45, 17).Value = '23/09/2023 18:00'
$ws.Cells.Item(45, 18).Value = 5.67
$ws.Cells.Item(45, 19).Value = '22/09/2023 15:12'
$ws.Cells.Item(45, 20).Value = 6.2
$ws.Cells.Item(45, 21).Value = '23/09/2023 18:00'
$ws.Cells.Item(45, 22).Value = 'https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/buducnost-petrovac/A7wXvB9R/'
